# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" positioned between "2021-Q3" and "总计",
#    populated with the per-fund holding breakdown for that quarter.
# 2. Insert a new summary row into the "总计" (totals) sheet for "2022-Q1".
#
# NOTE: worksheet/range handles in this host resolve by *position*, not by
# stable identity, so we always re-fetch sheets/ranges by name right after
# any operation that changes sheet count or row/column indices instead of
# reusing a variable captured beforehand.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q1" worksheet right after "2021-Q3"
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q3")
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$q1.Name = "2022-Q1"

# Headers (row 1) - copy the header formatting (bold + border, the style
# index already shared by the other sheets) from the totals sheet's header.
$wb.Worksheets.Item("总计").Range("B1").Copy()
$wb.Worksheets.Item("2022-Q1").Range("B1:H1").PasteSpecial(-4122)

$wb.Worksheets.Item("2022-Q1").Range("B1").Value = "基金代码"
$wb.Worksheets.Item("2022-Q1").Range("C1").Value = "基金名称"
$wb.Worksheets.Item("2022-Q1").Range("D1").Value = "基金规模"
$wb.Worksheets.Item("2022-Q1").Range("E1").Value = "股票总仓位"
$wb.Worksheets.Item("2022-Q1").Range("F1").Value = "仓位占比"
$wb.Worksheets.Item("2022-Q1").Range("G1").Value = "持有市值(亿元)"
$wb.Worksheets.Item("2022-Q1").Range("H1").Value = "仓位排名"

# Row 2 - fund 001914
$wb.Worksheets.Item("总计").Range("A2").Copy()
$wb.Worksheets.Item("2022-Q1").Range("A2").PasteSpecial(-4122)
$wb.Worksheets.Item("2022-Q1").Range("A2").Value = 0

$wb.Worksheets.Item("2022-Q1").Range("B2").Value = "'001914"
$wb.Worksheets.Item("2022-Q1").Range("B2").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("C2").Value = "中信建投聚利混合A"
$wb.Worksheets.Item("2022-Q1").Range("D2").Value = "'0.13"
$wb.Worksheets.Item("2022-Q1").Range("D2").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("E2").Value = "'39.07"
$wb.Worksheets.Item("2022-Q1").Range("E2").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("F2").Value = "'2.04"
$wb.Worksheets.Item("2022-Q1").Range("F2").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("G2").Value = "'0.0027"
$wb.Worksheets.Item("2022-Q1").Range("G2").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("H2").Value = 9

# Row 3 - fund 000041
$wb.Worksheets.Item("总计").Range("A3").Copy()
$wb.Worksheets.Item("2022-Q1").Range("A3").PasteSpecial(-4122)
$wb.Worksheets.Item("2022-Q1").Range("A3").Value = 1

$wb.Worksheets.Item("2022-Q1").Range("B3").Value = "'000041"
$wb.Worksheets.Item("2022-Q1").Range("B3").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("C3").Value = "华夏全球精选股票(QDII)"
$wb.Worksheets.Item("2022-Q1").Range("D3").Value = "'0.02"
$wb.Worksheets.Item("2022-Q1").Range("D3").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("E3").Value = "'39.07"
$wb.Worksheets.Item("2022-Q1").Range("E3").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("F3").Value = "'2.04"
$wb.Worksheets.Item("2022-Q1").Range("F3").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("G3").Value = "'0.0004"
$wb.Worksheets.Item("2022-Q1").Range("G3").Style = "Normal"
$wb.Worksheets.Item("2022-Q1").Range("H3").Value = 9

# ---------------------------------------------------------------------------
# Step 2: insert a new summary row for "2022-Q1" at the top of the "总计"
# sheet's data (row 2), pushing the existing quarters down. The row-index
# helper column (A) is renumbered sequentially (0, 1, 2, ...) afterwards.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("总计").Rows.Item(2).Insert()
$wb.Worksheets.Item("总计").Range("A2:D2").Style = "Normal"

$wb.Worksheets.Item("总计").Range("A3").Copy()
$wb.Worksheets.Item("总计").Range("A2").PasteSpecial(-4122)

$wb.Worksheets.Item("总计").Range("A2").Value = 0
$wb.Worksheets.Item("总计").Range("B2").Value = "2022-Q1"
$wb.Worksheets.Item("总计").Range("C2").Value = 2
$wb.Worksheets.Item("总计").Range("D2").Value = 0

$wb.Worksheets.Item("总计").Range("A3").Value = 1
$wb.Worksheets.Item("总计").Range("A4").Value = 2
